# Clear the "Student Name" / "Student ID" values that were filled in on
# the report header (row 2), and leave the selection on that range the
# way Excel does right after a Delete key-press.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:B2").ClearContents()

$ws.Range("A2:B2").Select()
